$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date value (2021-04-21, serial 44307) to cell C1, formatted as a
# short date (built-in numFmtId 14).
$ws.Range("C1").Value2 = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
